# Actualización automática 2025-06-20 15:15:09
#
# Registers a sale of 160.81 (LAVABOS group, advisor ALMEIDA CUATIN JHONATHANN
# CARLOS) that propagates across the three report sheets of the workbook:
#   - "VENTAS POR GRUPO"      : per-client/per-group sales matrix
#   - "VENTA MENSUAL"         : per-client/per-month sales matrix (junio)
#   - "CUMPLIMIENTO MENSUAL"  : per-group budget-vs-sales compliance summary

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ----------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Client row 8 (COMFALASDI COMPAÑIA FAMILIAR LASCANO DIAZ C. LTDA.), column I = LAVABOS
$wsGrupo.Range("I8").Value = 160.81

# Row 30 keeps a textual "<n> de 28" tally of non-zero clients per group;
# LAVABOS (column I) now has one more client with sales.
$wsGrupo.Range("I30").Value = "1 de 28"

# ----------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ----------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Same client row 8, column F = junio sales total.
$wsMensual.Range("F8").Value = 1990.27

# Row 30 totals column F (junio) across all clients.
$wsMensual.Range("F30").Value = 5792.25

# ----------------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL"
# ----------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# LAVABOS group row 8: VENTA, POR CUMPLIR and CUMPLIMIENTO (%) updated.
$wsCumpl.Range("D8").Value = 160.81
$wsCumpl.Range("E8").Value = 464.19
$wsCumpl.Range("F8").Value = 0.257296

# TOTAL row 19 aggregates VENTA, POR CUMPLIR and CUMPLIMIENTO (%).
$wsCumpl.Range("D19").Value = 5786.49
$wsCumpl.Range("E19").Value = 23751.30107555787
$wsCumpl.Range("F19").Value = 0.1959012434341525
